$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A23 holds a date-looking string ("06/08/2025"); force text formatting
# first so Excel doesn't auto-convert it to a date serial, then drop the
# format override back to Normal so no style is left behind on the cell.
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "06/08/2025"
$ws.Range("A23").Style = "Normal"

$ws.Range("B23").Value = "Grau"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = "Deportivo Garcilaso"
$ws.Range("F23").Value = "D"
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 1
$ws.Range("K23").Value = 0.89
$ws.Range("L23").Value = 1.65
$ws.Range("M23").Value = 16
$ws.Range("N23").Value = 16
$ws.Range("O23").Value = 5
$ws.Range("P23").Value = 5
